$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "POR"
$ws.Range("C2").Value = 13.63333333333333
$ws.Range("B3").Value = "NJN"
$ws.Range("C3").Value = 11.78571428571428
$ws.Range("B4").Value = "CLE"
$ws.Range("C4").Value = 12.43333333333334
$ws.Range("B5").Value = "DAL"
$ws.Range("C5").Value = 12.84166666666667
$ws.Range("B6").Value = "MIA"
$ws.Range("C6").Value = 10.5
$ws.Range("B7").Value = "SEA"
$ws.Range("C7").Value = 13.78125
$ws.Range("B8").Value = "ATL"
$ws.Range("C8").Value = 12.57058823529412
$ws.Range("B9").Value = "WAS"
$ws.Range("C9").Value = 14.36428571428571
$ws.Range("B10").Value = "MIL"
$ws.Range("C10").Value = 13.36153846153846
$ws.Range("B11").Value = "LAC"
$ws.Range("C11").Value = 12.76
$ws.Range("B12").Value = "SAS"
$ws.Range("C12").Value = 13.32666666666667
$ws.Range("B13").Value = "DET"
$ws.Range("C13").Value = 14.22666666666667
$ws.Range("B14").Value = "ORL"
$ws.Range("C14").Value = 12.54666666666667
$ws.Range("B15").Value = "UTA"
$ws.Range("C15").Value = 12.96428571428572
$ws.Range("B16").Value = "MEM"
$ws.Range("C16").Value = 10.43529411764706
$ws.Range("B17").Value = "HOU"
$ws.Range("C17").Value = 13.15714285714286
$ws.Range("B18").Value = "DEN"
$ws.Range("C18").Value = 12.92307692307692
$ws.Range("B19").Value = "LAL"
$ws.Range("C19").Value = 15.77142857142857
$ws.Range("B20").Value = "GSW"
$ws.Range("C20").Value = 13.25384615384615
$ws.Range("B21").Value = "IND"
$ws.Range("C21").Value = 13.08181818181818
$ws.Range("B22").Value = "CHI"
$ws.Range("C22").Value = 11.9
$ws.Range("B23").Value = "PHI"
$ws.Range("C23").Value = 10.93333333333333
$ws.Range("B24").Value = "CHH"
$ws.Range("C24").Value = 12.62857142857143
$ws.Range("B25").Value = "BOS"
$ws.Range("C25").Value = 12.05833333333334
$ws.Range("B26").Value = "TOR"
$ws.Range("C26").Value = 13.52
$ws.Range("B27").Value = "SAC"
$ws.Range("C27").Value = 14.60769230769231
$ws.Range("B28").Value = "PHO"
$ws.Range("C28").Value = 12.52727272727273
$ws.Range("B29").Value = "NYK"
$ws.Range("C29").Value = 11.02142857142857
$ws.Range("B30").Value = "MIN"
$ws.Range("C30").Value = 13.52142857142857
